# header changed adminpanel and login done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Admin panel task (row 3) is now marked as done.
$ws.Cells.Item(3, 2).Value = 1

# New tasks appended below the existing list (row 10 intentionally left blank).
$ws.Cells.Item(11, 1).Value = "تنظیمات تبلیغات"
$ws.Cells.Item(12, 1).Value = "تنظیم 3 پست بزرگ صفحه اصلی"
$ws.Cells.Item(13, 1).Value = "تنظیمات اسلایدر پایین صفحه"
$ws.Cells.Item(14, 1).Value = "افزودن کالا"
$ws.Cells.Item(15, 1).Value = "تنظیمات دسته بندی"
$ws.Cells.Item(16, 1).Value = "ادیت اطلاعات فروشگاه و توضیحات "
$ws.Cells.Item(17, 1).Value = "تغییر اطلاعات برای مدیر سایت"

$ws.Cells.Item(18, 1).Value = "صفحه ورود"
$ws.Cells.Item(18, 2).Value = 1

$ws.Cells.Item(19, 1).Value = "ثبت نام"
$ws.Cells.Item(20, 1).Value = "صفحه فراموشی کلمه عبور"

# Scroll the view down and select B20, matching the saved window state.
$ws.Range("B20").Select()
$excel.ActiveWindow.ScrollRow = 4
